# Updated cryptos list on Fri Nov  1 09:57:54 UTC 2024 with GitHub Actions
#
# Refresh Price (col D) / Volume(1h) (col E) figures for the crypto table,
# and swap the Bittensor / Fetch.AI rows (32/33) to reflect the new ranking
# order along with their refreshed figures.
#
# All of these source cells are plain text (inlineStr) in the workbook, not
# numbers - e.g. "69.523.38" / "0.0\u20830901" aren't valid numeric literals,
# and values like "153.20" or "578.68" must keep their exact text form
# (a bare numeric assignment would normalize/strip them, e.g. 153.20 -> 153.2).
# Forcing the cell to Text format before assigning keeps everything as text,
# and resetting the style back to Normal afterwards avoids leaving behind an
# extra/changed number-format style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "69.523.38"
Set-TextValue "E2" "  -3.88%  "
Set-TextValue "D3" "2.509.16"
Set-TextValue "E3" "  -4.83%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "578.68"
Set-TextValue "E5" "  -0.81%  "
Set-TextValue "E6" "  -4.80%  "
Set-TextValue "E7" "  +0.07%  "
Set-TextValue "D8" "0.523"
Set-TextValue "E8" "  +0.60%  "
Set-TextValue "D9" "2.508.35"
Set-TextValue "E9" "  -4.85%  "
Set-TextValue "D10" "0.158"
Set-TextValue "E10" "  -7.86%  "
Set-TextValue "E11" "  -1.42%  "
Set-TextValue "E12" "  -2.84%  "
Set-TextValue "E13" "  -1.15%  "
Set-TextValue "D14" "2.969.67"
Set-TextValue "E14" "  -4.85%  "
Set-TextValue "D15" "69.397.39"
Set-TextValue "E15" "  -3.85%  "
Set-TextValue "D17" "24.96"
Set-TextValue "E17" "  -3.34%  "
Set-TextValue "D18" "2.515.91"
Set-TextValue "E18" "  -4.87%  "
Set-TextValue "D19" "7.79"
Set-TextValue "E19" "  -1.45%  "
Set-TextValue "D20" "11.33"
Set-TextValue "E20" "  -6.30%  "
Set-TextValue "D21" "348.65"
Set-TextValue "E21" "  -7.09%  "
Set-TextValue "D22" "3.94"
Set-TextValue "E22" "  -4.05%  "
Set-TextValue "E23" "  -0.04%  "
Set-TextValue "D24" "1.93"
Set-TextValue "E24" "  -5.41%  "
Set-TextValue "D25" "68.71"
Set-TextValue "E25" "  -2.98%  "
Set-TextValue "D26" "3.98"
Set-TextValue "E26" "  -6.07%  "
Set-TextValue "D27" "8.94"
Set-TextValue "D28" "2.639.79"
Set-TextValue "E28" "  -4.85%  "
Set-TextValue "D29" "0.992"
Set-TextValue "E29" "  -0.70%  "
Set-TextValue "D30" "0.0₃0901"
Set-TextValue "E30" "  -4.85%  "
Set-TextValue "D31" "7.89"
Set-TextValue "E31" "  -0.80%  "
Set-TextValue "B32" "Fetch.AI"
Set-TextValue "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D32" "1.26"
Set-TextValue "E32" "  -0.23%  "
Set-TextValue "B33" "Bittensor"
Set-TextValue "C33" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D33" "465.11"
Set-TextValue "E33" "  -6.04%  "
Set-TextValue "E34" "  -2.05%  "
Set-TextValue "E35" "  -0.02%  "
Set-TextValue "D36" "0.118"
Set-TextValue "E36" "  +3.85%  "
Set-TextValue "D37" "153.20"
Set-TextValue "E37" "  -5.80%  "
Set-TextValue "D39" "18.37"
Set-TextValue "E39" "  -4.15%  "
Set-TextValue "D41" "4.75"
Set-TextValue "E41" "  -2.82%  "
Set-TextValue "E42" "  -2.38%  "
Set-TextValue "E43" "  -7.10%  "
Set-TextValue "D44" "1.16"
Set-TextValue "E44" "  -14.36%  "
Set-TextValue "E45" "  -10.40%  "
Set-TextValue "D46" "38.13"
Set-TextValue "E46" "  -2.25%  "
Set-TextValue "D47" "143.13"
Set-TextValue "E47" "  -5.45%  "
Set-TextValue "D48" "0.529"
Set-TextValue "E48" "  -2.58%  "
Set-TextValue "D49" "3.49"
Set-TextValue "E49" "  -4.22%  "
Set-TextValue "E50" "  -4.43%  "
Set-TextValue "D51" "0.0732"
Set-TextValue "E51" "  -1.98%  "
